# A new daily price record was inserted right before the existing row 115
# (the Fecha/D column holds the date serial 44574), pushing every
# subsequent record down by one row (old row 115 -> new row 116, ...,
# old row 198 -> new row 199). The new row keeps the same Volumen /
# Precio / Unidad / Origen / Precio $/Kg / Kg o Unidades / Clasificacion
# values as the record that used to sit at row 115 - only the date
# changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 115:198 down to 116:199, leaving a blank row 115 behind.
$ws.Rows.Item(115).Insert()

# Duplicate the record that is now at row 116 (the original row-115
# record) into the freshly inserted row 115.
$src = $ws.Range("A116:R116")
$dst = $ws.Range("A115:R115")
$src.Copy($dst)

# The new record carries a different date (44574) than the one it was
# copied from.
$ws.Range("D115").Value = 44574
